# Apply the edits to rows 160-178 (columns C and D) as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D-column value updates (rows 160-178) ---
$ws.Range("D160").Value = 88.784999999999997
$ws.Range("D161").Value = 89.17
$ws.Range("D162").Value = 89.605000000000004
$ws.Range("D163").Value = 90.05
$ws.Range("D164").Value = 90.56
$ws.Range("D165").Value = 90.795000000000002
$ws.Range("D166").Value = 91.015000000000001
$ws.Range("D167").Value = 91.564999999999998
$ws.Range("D168").Value = 91.614999999999995
$ws.Range("D169").Value = 92.125
$ws.Range("D170").Value = 92.44
$ws.Range("D171").Value = 92.75
$ws.Range("D172").Value = 93.155000000000001
$ws.Range("D173").Value = 93.155000000000001
$ws.Range("D174").Value = 93.155000000000001
$ws.Range("D175").Value = 93.155000000000001
$ws.Range("D176").Value = 101.22
$ws.Range("D177").Value = 101.545
$ws.Range("D178").Value = 102.83499999999999

# --- C-column: turn these into formulas that reference the D cell one row up ---
# C161 and C177 are plain (non-shared) formulas; C162-C176 and C178 form the
# "fill-down" shared-formula run (C162:C178) referencing the D cell directly above.
$ws.Range("C161").Formula = "=D160"
$ws.Range("C162:C178").Formula = "=D161"
$ws.Range("C177").Formula = "=D176"

# --- Clear the bold/red style that used to be applied to C172:C178 ---
# (B/D/E/F on those rows keep their style; only column C reverts to Normal.)
$ws.Range("C172:C178").Style = "Normal"

# --- Leftover view state (zoom + selection) from the author's last save ---
$ws.Range("E7").Select() | Out-Null
$excel.ActiveWindow.Zoom = 159
